$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Termiet" week block (rows 74-80): fill in the previously-blank
# "Maandag" attendance row (75) with 4 hours for every person, matching the
# already-filled "Vrijdag" row (79) style/template.
$ws.Range("C14:H14").Copy()
$ws.Range("C75:H75").PasteSpecial(-4122)
$ws.Range("C75:H75").Value = 4

# Weekly totals row (80) — Totaal Game-Lab uren is a hand-typed literal,
# not a formula, so it needs to be bumped explicitly; C80:H80 are
# SUM(...) formulas over C75:C79 etc. and recalc on their own.
$ws.Range("B80").Value = 21

# --- Same fix for the newest week block (rows 162-168): row 167 ("Vrijdag")
# was left blank; fill it in with 4 hours each, copying the style from the
# already-filled row 79 template (thick-bottom border variant).
$ws.Range("C79:H79").Copy()
$ws.Range("C167:H167").PasteSpecial(-4122)
$ws.Range("C167:H167").Value = 4

# Weekly totals row (168) — same literal-total fix as row 80.
$ws.Range("B168").Value = 23

# Selection left where the author was last working.
$null = $ws.Range("L22").Select()
